$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 1 text updates: "address book" -> "source manager" wording.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(1)

# "TextBox 47": "[command commits address book]" -> "[command commits source manager]"
$tb47 = $s.Shapes.Item("TextBox 47")
$tb47Height = $tb47.Height
$tr47 = $tb47.TextFrame.TextRange
# Replace only the second run's text ("command commits address book]"),
# leaving the leading "[" run untouched so formatting/structure survives.
$run2 = $tr47.Characters(2, $tr47.Text.Length - 1)
$run2.Text = "command commits source manager]"
# The autofit text box keeps the same rendered height in the source deck;
# restore it so the shape extent isn't nudged by a rounding-level re-measure.
$tb47.Height = $tb47Height

# "Rounded Rectangle 50": "Purge redundant states and then save address book to addressBookStateList "
#   -> "Purge redundant states and then save source manager to sourceManagerStateList "
$rr50 = $s.Shapes.Item("Rounded Rectangle 50")
$tr50 = $rr50.TextFrame.TextRange
# Edit back-to-front so earlier character offsets stay valid.
$runB = $tr50.Characters(54, 20)
$runB.Text = "sourceManagerStateList"
$runA = $tr50.Characters(1, 53)
$runA.Text = "Purge redundant states and then save source manager to "

# ---------------------------------------------------------------------------
# 2) Refresh the cached "datetimeFigureOut" text on every Date Placeholder
#    (slide master + all slide layouts) from 6/7/2018 to 15/4/19.
# ---------------------------------------------------------------------------
function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            $len = $tr.Text.Length
            if ($len -gt 0) {
                $chars = $tr.Characters(1, $len)
                $chars.Text = "15/4/19"
            }
        }
    }
}

Update-DatePlaceholders $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    Update-DatePlaceholders $layouts.Item($l).Shapes
}
